# Insert a new weekly data row at row 192, shifting the existing rows
# 192-300 down to 193-301 (dimension grows from A1:R300 to A1:R301).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 192; this pushes rows 192..300 to 193..301.
$ws.Rows(192).Insert()

# Populate the newly inserted row 192 with the new record's data.
$ws.Cells.Item(192, 1).Value = 9
$ws.Cells.Item(192, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(192, 3).Value = "Metropolitana"
$ws.Cells.Item(192, 4).Value = Get-Date -Year 2021 -Month 11 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(192, 5).Value = 13
$ws.Cells.Item(192, 6).Value = 100112039
$ws.Cells.Item(192, 7).Value = "Ciboulette"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 160
$ws.Cells.Item(192, 11).Value = 1000
$ws.Cells.Item(192, 12).Value = 1200
$ws.Cells.Item(192, 13).Value = 1100
$ws.Cells.Item(192, 14).Value = "`$/docena de atados"
$ws.Cells.Item(192, 15).Value = "Región Metropolitana"
$ws.Cells.Item(192, 16).Value = 367
$ws.Cells.Item(192, 17).Value = 3
$ws.Cells.Item(192, 18).Value = "Hortaliza"
